$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 20 -----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Klacht over levering"
$logs.Range("B20").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C20").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D20").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E20").Value = "Bedankt, we hebben dit doorgestuurd naar klachten@testbedrijf123.nl."
$logs.Range("F20").Value = "2025-08-14 21:03:37"
$logs.Range("G20").Value = "Nee"
$logs.Range("H20").Value = "Ja"
$logs.Range("I20").Value = "Nee"
$logs.Range("J20").Value = "Nee"

# --- extend conditional formatting ranges to cover the new row -------------
$dFc = $logs.Range("D2:D19").FormatConditions
$dFc.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))

$gFc = $logs.Range("G2:G19").FormatConditions
$gFc.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))

$hFc = $logs.Range("H2:H19").FormatConditions
$hFc.Item(1).ModifyAppliesToRange($logs.Range("H2:H20"))

$iFc = $logs.Range("I2:I19").FormatConditions
$iFc.Item(1).ModifyAppliesToRange($logs.Range("I2:I20"))

$jFc = $logs.Range("J2:J19").FormatConditions
$jFc.Item(1).ModifyAppliesToRange($logs.Range("J2:J20"))

# --- Dashboard sheet: bump count for "Intern verzoek / Actie voor medewerker"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 14
